$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $TextValue)
    $range = $ws.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $TextValue
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "37.128.63"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.047.74"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "249.65"
$ws.Range("E5").Value = "  -3.09%  "
Set-TextValue "D6" "0.656"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  -0.12%  "
Set-TextValue "D8" "55.14"
$ws.Range("E8").Value = "  +16.55%  "
Set-TextValue "D9" "61.95"
$ws.Range("E9").Value = "  +0.47%  "
Set-TextValue "D10" "0.379"
$ws.Range("E10").Value = "  +0.84%  "
Set-TextValue "D11" "0.0756"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  +5.46%  "
Set-TextValue "D13" "15.15"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "2.343.82"
$ws.Range("E14").Value = "  -3.59%  "
Set-TextValue "D15" "0.826"
$ws.Range("E15").Value = "  -3.63%  "
Set-TextValue "D16" "5.27"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "2.050.48"
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("D18").Value = "37.027.79"
$ws.Range("E18").Value = "  +0.27%  "
Set-TextValue "D19" "72.47"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").Value = "0.0₃0879"
$ws.Range("E20").Value = "  +3.63%  "
Set-TextValue "D21" "14.42"
$ws.Range("E21").Value = "  +6.59%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D22" "5.29"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D23" "238.24"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -1.68%  "
Set-TextValue "D26" "170.34"
$ws.Range("E26").Value = "  -1.48%  "
Set-TextValue "D27" "9.17"
$ws.Range("E27").Value = "  -1.31%  "
Set-TextValue "D28" "20.33"
$ws.Range("E28").Value = "  -5.45%  "
Set-TextValue "D29" "2.01"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  -1.01%  "
Set-TextValue "D31" "4.61"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  +15.36%  "
Set-TextValue "D33" "0.0628"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "4.40"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("B35").Value = "Gas"
$ws.Range("C35").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextValue "D35" "19.61"
$ws.Range("E35").Value = "  -16.54%  "
$ws.Range("E36").Value = "  +0.14%  "
Set-TextValue "D37" "2.29"
$ws.Range("E37").Value = "  -5.57%  "
Set-TextValue "D38" "0.0838"
$ws.Range("E38").Value = "  -12.55%  "
$ws.Range("E39").Value = "  -5.70%  "
Set-TextValue "D40" "0.112"
$ws.Range("E40").Value = "  +32.53%  "
Set-TextValue "D41" "1.35"
$ws.Range("E41").Value = "  -1.77%  "
Set-TextValue "D42" "18.14"
$ws.Range("E42").Value = "  +10.73%  "
$ws.Range("E43").Value = "  +0.22%  "
Set-TextValue "D44" "1.14"
$ws.Range("E44").Value = "  -4.74%  "
Set-TextValue "D45" "97.43"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  +63.11%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "1.304.77"
$ws.Range("E48").Value = "  -4.52%  "
Set-TextValue "D49" "2.38"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  +2.47%  "
Set-TextValue "D51" "6.89"
$ws.Range("E51").Value = "  -4.40%  "

Write-Host "Applied all changes"
